# amiris_data_structure.xlsx - "reorganizing amiris config files"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) load_shedding sheet: VOLL values replaced by simpler numbers, and the
#    TimeSeries (csv-path) column is replaced by plain 0s.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("load_shedding")
$ws.Range("B2").Value = 47
$ws.Range("C2").Value = 0
$ws.Range("B3").Value = 4000
$ws.Range("C3").Value = 0
$ws.Range("B4").Value = 1500
$ws.Range("C4").Value = 0
$ws.Range("B5").Value = 250
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = 500
$ws.Range("C6").Value = 0

# ---------------------------------------------------------------------------
# 2) scenario_data_emlab sheet: move the OTHER column from the end (L) to
#    right after CO2 (new column D), shifting the fuel columns right by one.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("scenario_data_emlab")
$ws.Range("D1").Value = "OTHER"
$ws.Range("E1").Value = "HARD_COAL"
$ws.Range("F1").Value = "OIL"
$ws.Range("G1").Value = "HYDROGEN"
$ws.Range("H1").Value = "LIGNITE"
$ws.Range("I1").Value = "NATURAL_GAS"
$ws.Range("J1").Value = "NUCLEAR"
$ws.Range("K1").Value = "WASTE"
$ws.Range("L1").Value = "BIOMASS"

$ws.Range("C2").Value = 167.9999999999998

$ws.Range("D3").Value = 64.47999999998137
$ws.Range("E3").Value = 6.730000000001382
$ws.Range("F3").Value = 79.69000000000142
$ws.Range("G3").Value = 48.57925285126839
$ws.Range("H3").Value = 6.479999999999997
$ws.Range("I3").Value = 14.65000000003783
$ws.Range("J3").Value = 1.689999999999999
$ws.Range("K3").Value = 7.499999999999996
$ws.Range("L3").Value = 34.99999999995634

# ---------------------------------------------------------------------------
# 3) conventionals sheet: row 2 becomes an OTHER plant, and 8 additional
#    identical OTHER rows are appended.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("conventionals")
$ws.Range("B2").Value = 20383300053
$ws.Range("C2").Value = "OTHER"
$ws.Range("D2").Value = 1.5
$ws.Range("E2").Value = 0.4
$ws.Range("F2").Value = 1500
$ws.Range("G2").Value = 1500

$otherRows = @(
    @(1, 20323300052),
    @(2, 20313300051),
    @(3, 20303300050),
    @(4, 20293300049),
    @(5, 20283300048),
    @(6, 20273300047),
    @(7, 20263300046),
    @(8, 20253300045)
)
foreach ($r in $otherRows) {
    $idx = $r[0]
    $identifier = $r[1]
    $row = 2 + $idx
    $ws.Cells.Item($row, 1).Value = $idx
    $ws.Cells.Item($row, 2).Value = $identifier
    $ws.Cells.Item($row, 3).Value = "OTHER"
    $ws.Cells.Item($row, 4).Value = 1.5
    $ws.Cells.Item($row, 5).Value = 0.4
    $ws.Cells.Item($row, 6).Value = 1500
    $ws.Cells.Item($row, 7).Value = 1500
}

# ---------------------------------------------------------------------------
# 4) renewables sheet: update identifiers/opex/set for all existing rows and
#    append three additional WindOff rows at the bottom.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("renewables")
$renewables = @(
    @(0,  20390200023, 4000, 1.35, "WindOn"),
    @(1,  20390300033, 5000, 0,    "OtherPV"),
    @(2,  20380300032, 5000, 0,    "OtherPV"),
    @(3,  20370300031, 5000, 0,    "OtherPV"),
    @(4,  20360300030, 5500, 0,    "OtherPV"),
    @(5,  20350300029, 4500, 0,    "OtherPV"),
    @(6,  20340200022, 4000, 1.35, "WindOn"),
    @(7,  20340300028, 3257, 0,    "OtherPV"),
    @(8,  20340100043, 7000, 2.7,  "WindOff"),
    @(9,  20330300027, 4500, 0,    "OtherPV"),
    @(10, 20330100042, 7000, 2.7,  "WindOff"),
    @(11, 20320300026, 4500, 0,    "OtherPV"),
    @(12, 20320100041, 7000, 2.7,  "WindOff"),
    @(13, 20310300025, 4500, 0,    "OtherPV"),
    @(14, 20310100040, 7000, 2.7,  "WindOff"),
    @(15, 20300200021, 4000, 1.35, "WindOn"),
    @(16, 20300300024, 4500, 0,    "OtherPV"),
    @(17, 20300100039, 7000, 2.7,  "WindOff"),
    @(18, 20290100038, 7000, 2.7,  "WindOff"),
    @(19, 20280100037, 7000, 2.7,  "WindOff"),
    @(20, 20270100036, 7000, 2.7,  "WindOff"),
    @(21, 20260100035, 7000, 2.7,  "WindOff"),
    @(22, 20250100034, 7000, 2.7,  "WindOff")
)
foreach ($r in $renewables) {
    $idx = $r[0]
    $row = 2 + $idx
    $ws.Cells.Item($row, 1).Value = $idx
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = "NONE"
    $ws.Cells.Item($row, 7).Value = "-"
    $ws.Cells.Item($row, 8).Value = "-"
    $ws.Cells.Item($row, 9).Value = "-"
}

# ---------------------------------------------------------------------------
# 5) storages sheet: the single remaining storage keeps the old row 2 slot
#    but gets new data, and the second storage row is deleted entirely.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("storages")
$ws.Range("B2").Value = 20352600044
$ws.Range("C2").Value = "STORAGE"
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 0.92
$ws.Range("F2").Value = 0.92
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1700
$ws.Range("A3:H3").Delete()
